$wb = $excel.ActiveWorkbook

# Update timestamps on the "zh-cn" handback status sheet
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-24 03:19:22"
$wsZhCn.Range("H2").Value = "2016-03-24 03:20:05"

# Update timestamps on the "de-de" handback status sheet
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-24 03:19:26"
$wsDeDe.Range("H2").Value = "2016-03-24 03:20:12"
